$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Final feed: 7 articles (was 12). Two brand-new headlines land on top
# (A1:A2), five headlines that were already in the feed move up to
# A3:A7, and the five oldest articles (old rows 8-12) drop off the
# bottom entirely.
# ---------------------------------------------------------------------------
$newText = @(
    "Екипажът на Союз МС-22 ще остане още няколко месеца на МКС",
    "„Първата извънземна сонда, достигнала Земята, ще бъде твърде сложна, за да я разберем“",
    "Заглушават GPS сигналите покрай българското крайбрежие – източникът за сега е неизвестен",
    "Lenovo ThinkPad Z13: Еволюция в действие",
    "НАСА работи върху хибридната мисия до Титан и още над дузина проекти в космоса",
    "Историческото първо изстрелване на ракета от британска територия претърпя провал",
    "HWO е новият космически телескоп за търсене на извънземен живот"
)

$newUrl = @(
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/%d0%ba%d0%be%d1%81%d0%bc%d0%be%d1%81/%d0%b5%d0%ba%d0%b8%d0%bf%d0%b0%d0%b6%d1%8a%d1%82-%d0%bd%d0%b0-%d1%81%d0%be%d1%8e%d0%b7-%d0%bc%d1%81-22-%d1%89%d0%b5-%d0%be%d1%81%d1%82%d0%b0%d0%bd%d0%b5-%d0%be%d1%89%d0%b5-%d0%bd%d1%8f%d0%ba%d0%be%d0%bb%d0%ba%d0%be-%d0%bc%d0%b5%d1%81%d0%b5%d1%86%d0%b0-%d0%bd%d0%b0-%d0%bc%d0%ba%d1%81-404300.html",
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/%d0%ba%d0%be%d1%81%d0%bc%d0%be%d1%81/%d0%bf%d1%8a%d1%80%d0%b2%d0%b0%d1%82%d0%b0-%d0%b8%d0%b7%d0%b2%d1%8a%d0%bd%d0%b7%d0%b5%d0%bc%d0%bd%d0%b0-%d1%81%d0%be%d0%bd%d0%b4%d0%b0-%d0%b4%d0%be%d1%81%d1%82%d0%b8%d0%b3%d0%bd%d0%b0%d0%bb%d0%b0-%d0%b7%d0%b5%d0%bc%d1%8f%d1%82%d0%b0-%d1%89%d0%b5-%d0%b1%d1%8a%d0%b4%d0%b5-%d1%82%d0%b2%d1%8a%d1%80%d0%b4%d0%b5-404299.html",
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/%d0%b7%d0%b0%d0%b3%d0%bb%d1%83%d1%88%d0%b0%d0%b2%d0%b0%d1%82-gps-%d1%81%d0%b8%d0%b3%d0%bd%d0%b0%d0%bb%d0%b8%d1%82%d0%b5-%d0%bf%d0%be%d0%ba%d1%80%d0%b0%d0%b9-%d0%b1%d1%8a%d0%bb%d0%b3%d0%b0%d1%80%d1%81-403838.html",
    "https://www.kaldata.com/%d1%80%d0%b5%d0%b2%d1%8e%d1%82%d0%b0/lenovo-thinkpad-z13-%d0%b5%d0%b2%d0%be%d0%bb%d1%8e%d1%86%d0%b8%d1%8f-%d0%b2-%d0%b4%d0%b5%d0%b9%d1%81%d1%82%d0%b2%d0%b8%d0%b5-399936.html",
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/%d0%ba%d0%be%d1%81%d0%bc%d0%be%d1%81/%d0%bd%d0%b0%d1%81%d0%b0-%d1%80%d0%b0%d0%b1%d0%be%d1%82%d0%b8-%d0%b2%d1%8a%d1%80%d1%85%d1%83-%d1%85%d0%b8%d0%b1%d1%80%d0%b8%d0%b4%d0%bd%d0%b0%d1%82%d0%b0-%d0%bc%d0%b8%d1%81%d0%b8%d1%8f-%d0%b4%d0%be-404112.html",
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/%d0%b8%d1%81%d1%82%d0%be%d1%80%d0%b8%d1%87%d0%b5%d1%81%d0%ba%d0%be%d1%82%d0%be-%d0%bf%d1%8a%d1%80%d0%b2%d0%be-%d0%b8%d0%b7%d1%81%d1%82%d1%80%d0%b5%d0%bb%d0%b2%d0%b0%d0%bd%d0%b5-%d0%bd%d0%b0-%d1%80-403804.html",
    "https://www.kaldata.com/it-%d0%bd%d0%be%d0%b2%d0%b8%d0%bd%d0%b8/hwo-%d0%b5-%d0%bd%d0%be%d0%b2%d0%b8%d1%8f%d1%82-%d0%ba%d0%be%d1%81%d0%bc%d0%b8%d1%87%d0%b5%d1%81%d0%ba%d0%b8-%d1%82%d0%b5%d0%bb%d0%b5%d1%81%d0%ba%d0%be%d0%bf-%d0%b7%d0%b0-%d1%82%d1%8a%d1%80%d1%81-404262.html"
)

$tooltip = "open this article"
$keepCount = 7

# ---------------------------------------------------------------------------
# 1. Write the new headline text into A1:A7 (shared strings get rebuilt from
#    these values automatically on save).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $keepCount; $i++) {
    $ws.Range("A" + ($i + 1)).Value = $newText[$i]
}

# Capture the "Normal"-ish per-cell style used throughout the list (style
# index 1 in the original file) before we touch the hyperlinks, so we can
# re-apply it exactly once the links are rebuilt.
$styleSample = $ws.Range("A2").Style

# ---------------------------------------------------------------------------
# 2. Drop the five oldest articles (rows 8-12): shrinks sheetData + the
#    dimension down to A1:A7.
# ---------------------------------------------------------------------------
$ws.Rows(($keepCount + 1).ToString() + ":12").Delete()

# ---------------------------------------------------------------------------
# 3. Rebuild the hyperlinks collection. The engine does not support removing
#    individual hyperlinks, so clear the (now stale) collection entirely and
#    re-add only the ones that belong to the articles that remain,
#    restoring the original cell style afterwards (Add() otherwise stamps
#    its own hyperlink style onto the cell).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
for ($i = 0; $i -lt $keepCount; $i++) {
    $cell = $ws.Range("A" + ($i + 1))
    [void]$ws.Hyperlinks.Add($cell, $newUrl[$i], [Type]::Missing, $tooltip)
    $cell.Style = $styleSample
}

Write-Output "done"
